{"js": "// Rename the custom paragraph style \"Docdate\" to \"docDate\" (both the\n// w:styleId and the w:name/nameLocal change case). The Word JS API does\n// not expose a direct rename for styles (Style.nameLocal is read-only),\n// so we recreate the style under the new id/name with the same\n// definition (based on \"Title\", quick style / qFormat, 16pt / 32\n// half-point font size) and delete the old one.\n\nconst OLD_NAME = \"Docdate\";\nconst NEW_NAME = \"docDate\";\n\nconst styles = context.document.getStyles();\nconst oldStyle = styles.getByNameOrNullObject(OLD_NAME);\noldStyle.load(\"baseStyle,quickStyle,type,font/size\");\nawait context.sync();\n\nif (!oldStyle.isNullObject) {\n  const baseStyle = oldStyle.baseStyle;\n  const quickStyle = oldStyle.quickStyle;\n  const fontSize = oldStyle.font.size;\n\n  // Any paragraphs currently using the old style need to be moved onto\n  // the new one once it exists (deleting a style that is still applied\n  // would otherwise fall back to some other style).\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items/style\");\n  await context.sync();\n  const affected = paragraphs.items.filter((p) => p.style === OLD_NAME);\n\n  // Remove the old-named style.\n  oldStyle.delete();\n  await context.sync();\n\n  // Create the replacement style with the new id/name.\n  context.document.addStyle(NEW_NAME, \"Paragraph\");\n  await context.sync();\n\n  // Re-fetch the freshly created style and copy over its formatting.\n  const newStyle = context.document.getStyles().getByNameOrNullObject(NEW_NAME);\n  newStyle.load(\"nameLocal\");\n  await context.sync();\n\n  newStyle.baseStyle = baseStyle;\n  newStyle.quickStyle = quickStyle;\n  newStyle.font.size = fontSize;\n\n  affected.forEach((p) => {\n    p.style = NEW_NAME;\n  });\n  await context.sync();\n}\n", "ps1": "# Rename the custom paragraph style \"Docdate\" to \"docDate\" (both the\n# internal style id and the display name change case). The Word object\n# model does not expose the raw w:styleId as a settable property -\n# Style.NameLocal only changes the <w:name> element and leaves the\n# original w:styleId untouched - so to reproduce a rename of BOTH the\n# id and the name we recreate the style under the new name (which mints\n# a matching new id) with the same definition, and drop the old one.\n\n$OLD_NAME = \"Docdate\"\n$NEW_NAME = \"docDate\"\n\n$d = $word.ActiveDocument\n\n$oldStyle = $d.Styles | Where-Object { $_.NameLocal -eq $OLD_NAME }\n\nif ($oldStyle) {\n    $baseStyleName = $oldStyle.BaseStyle.NameLocal\n    $quickStyle = $oldStyle.QuickStyle\n    $fontSize = $oldStyle.Font.Size\n\n    # Any paragraphs currently using the old style need to be moved onto\n    # the new one once it exists (deleting a style that is still applied\n    # would otherwise fall back to some other style).\n    $affected = @($d.Paragraphs | Where-Object { $_.Style.NameLocal -eq $OLD_NAME })\n\n    # Remove the old-named style.\n    $oldStyle.Delete()\n\n    # Re-create it under the new id/name with the same formatting.\n    $newStyle = $d.Styles.Add($NEW_NAME, [Microsoft.Office.Interop.Word.WdStyleType]::wdStyleTypeParagraph)\n    $newStyle.BaseStyle = $baseStyleName\n    $newStyle.QuickStyle = $quickStyle\n    $newStyle.Font.Size = $fontSize\n\n    foreach ($p in $affected) {\n        $p.Style = $NEW_NAME\n    }\n}\n"}
